$p = $ppt.ActivePresentation
$s = $p.Slides.Item(37)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# 1. Shrink the autofit font scale (77.5% -> 70%) on the body placeholder.
#    PowerPoint normally computes this internally while the user types with
#    "Shrink text on overflow" active; best-effort set via the documented
#    COM property in case the host recomputes layout.
$tf.AutofitFontScale = 0.7
$tf.LineSpaceReduction = 0.2

# 2. Reword the "Resolve using ‘Mine’" bullet.
$mine = $tr.Find(" will discard whatever changes were in the other commit/branch.")
$mine.Text = " will discard whatever changes were in the commit/branch you’re merging in."

# 3. Reword the "Resolve Using ‘Theirs’" bullet.
$theirs = $tr.Find(" will discard whatever changes were made in the commit/branch you're currently on.")
$theirs.Text = " will discard whatever changes were in the commit/branch you’re merging into."

# 4. Remove the extra paragraph spacing after the last bullet ("If it’s a Unity file...").
$lastPara = $tr.Paragraphs(5,1)
$lastPara.ParagraphFormat.SpaceAfter = 0
